# Case_2_15 380 kV update: rewrite the B:O numeric block (rows 2-25) on the
# active sheet (pl_mw results grid) with the newly recomputed power-flow values.
# Columns D, K, N stay 0 for every row; column A (the index 0..23) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","E","F","G","H","I","J","L","M","O")

$newValues = @(
    @(1.003893904969004,0.3801705754779334,0.2848823815215447,1.569863645773403,0.4109121304418508,0.5802671699172706,0.6098850551920094,0.04392632853183187,0.5221184818567934,0.3101652500701704,1.930216354138139),
    @(0.8957952509933307,0.3712138418170952,0.2869908685361739,1.578743995872692,0.4154041102034753,0.5863010352372555,0.6199450980176167,0.04068430197565931,0.5152224545495585,0.2898550218144607,1.952119864615085),
    @(0.8292529174645438,0.3657349242427159,0.2883831781299016,1.585106861799048,0.4185275613623958,0.5903062601096636,0.6265200045280501,0.03868467560329236,0.5111518549782375,0.2774061664784071,1.966963508082799),
    @(0.8020958518027328,0.3635076344097143,0.2889751554908155,1.587928760011735,0.4198921134052682,0.5920139720025759,0.6292993353251788,0.03786760143702139,0.5095344125854808,0.2723390213553074,1.973362841341554),
    @(0.7975840548181736,0.3631381288480213,0.289074939680285,1.588411166075147,0.420124231448284,0.5923020996153809,0.6297668787526121,0.03773179503234303,0.50926834215349,0.2714979919571903,1.974446601745811),
    @(0.8288868291306812,0.3657048639858544,0.2883910620929235,1.585143991676553,0.4185455930182727,0.5903289849961553,0.6265570827493363,0.03867366513092207,0.5111298738209911,0.2773378049022455,1.967048393266566),
    @(0.9666578271821891,0.3770781715373062,0.285589133340137,1.572736703735877,0.4123850439541386,0.5822853036974394,0.6132711078837794,0.04281037967830059,0.5197069100622258,0.3031580147084796,1.937479010844541),
    @(1.235405829526428,0.399534705045653,0.2808681208265753,1.55562644898297,0.403210164803923,0.5688947882055118,0.5903784301529491,0.05084893549559411,0.5378158333684127,0.3539483421882892,1.890574734253647),
    @(1.431904927450887,0.4161160228802316,0.2778690273095989,1.547455612898119,0.3982508064128254,0.5605089222095856,0.5754900861810714,0.05670784268330209,0.5518967115801985,0.3913429809836799,1.862889483020282),
    @(1.52107468427846,0.4236749972669998,0.2766061400218902,1.544693658175476,0.3963836083495309,0.5570091600199376,0.5691370719433273,0.05936257184144722,0.5584690658169222,0.4083686354342149,1.851770239110067),
    @(1.554807716247296,0.4265394732096865,0.2761424659535923,1.543785058128663,0.395732631166382,0.5557291837157976,0.5667917948451819,0.06036628796277199,0.5609816430900878,0.4148175574747981,1.847772111387854),
    @(1.547544233986343,0.4259224695856574,0.2762416796386251,1.543974635884609,0.3958703333449662,0.5560028342927126,0.5672942010070834,0.06015019061175053,0.560439461118861,0.4134285989985784,1.848623722951686),
    @(1.523850607320412,0.4239106196463354,0.2765677017631454,1.544616155772417,0.3963289269335988,0.5569029475114888,0.5689429120208178,0.05944518002192467,0.5586753020011201,0.4088991609924193,1.851437049884254),
    @(1.509333134244287,0.4226785638443857,0.2767692942072095,1.545026983641719,0.3966171385322568,0.5574601935270564,0.5699606735372349,0.05901313416038789,0.5575977922717357,0.4061249550004717,1.853187977886947),
    @(1.426072887682153,0.4156223275339244,0.2779535981799679,1.547655327181872,0.3983806716826948,0.5607439795353031,0.5759137267693468,0.0565341331372835,0.5514705314376442,0.3902305723407906,1.863645864907426),
    @(1.374937864656431,0.4112974861371015,0.2787060842290963,1.549512304535412,0.3995622586265881,0.5628391644579267,0.5796733227489135,0.0550106099669847,0.5477542369755497,0.3804833279881592,1.870439485268463),
    @(1.345505898918134,0.4088114773245195,0.2791484417568011,1.550670280137219,0.4002784651516649,0.5640739091079112,0.5818752469680284,0.05413333337814663,0.5456324450572509,0.374878374918211,1.874485798544086),
    @(1.335537294732603,0.4079700271755371,0.2792998570986569,1.551077790751854,0.4005272391871486,0.5644970645010687,0.5826275621941548,0.05383613501822992,0.5449167501042069,0.3729808883997308,1.875879637988305),
    @(1.380383407759552,0.4117577168267132,0.2786249928927944,1.549305322959952,0.3994326885806174,0.5626130597433203,0.5792690180239131,0.05517289412749449,0.5481482173242256,0.381520796593172,1.86970192434913),
    @(1.530810929511006,0.4245014949255221,0.2764715464070591,1.544424000052004,0.396192703216677,0.5566373327267442,0.5684570035816936,0.05965230186673409,0.5591928350307569,0.4102295246531042,1.850604937538947),
    @(1.628927049810727,0.4328421616873186,0.2751489625164858,1.542033999718981,0.3944022024190517,0.5529959418515702,0.5617432518383421,0.06257066219112772,0.5665495710540682,0.4290019103017357,1.83936268074109),
    @(1.576579367002296,0.4283895871903383,0.2758470993059099,1.543236379405869,0.3953278467214574,0.5549152531180965,0.5652942154371008,0.06101393910984143,0.5626105510266797,0.4189820003809785,1.845249410741147),
    @(1.377921582715771,0.4115496452663194,0.2786616239345214,1.549398617819236,0.3994911522820743,0.5627151876300402,0.5794516780984331,0.05509952976115784,0.547970052831019,0.3810517605319177,1.870034937804064),
    @(1.162863196621402,0.393444337661947,0.2820626872505585,1.559482399571706,0.4053800843863655,0.5722622168284204,0.5962326940192924,0.04868239270247443,0.5378158333684127,0.3401931438713603,1.902075193318709)
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
